$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "Test Result" column header in column L, row 1
$ws.Range("L1").Value = "Test Result"

# Apply the built-in "Bad" cell style (pink fill / dark red text) to the new header cell
$ws.Range("L1").Style = "Bad"

# Match the font size used by the rest of the header row (16pt)
$ws.Range("L1").Font.Size = 16

# Reflect the post-edit selection: user clicked on L7 after entering the new header
[void]$ws.Range("L7").Select()
